$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D stays formatted as Text so numeric-looking strings
# (e.g. "29.647.35", "22.60", "0.00001062") are not reinterpreted as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.647.35"
$ws.Range("E2").Value = "  +8.39%  "
$ws.Range("D3").Value = "1.943.68"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "341.36"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D7").Value = "0.4771"
$ws.Range("E7").Value = "  +4.77%  "
$ws.Range("D8").Value = "0.4131"
$ws.Range("E8").Value = "  +8.73%  "
$ws.Range("D9").Value = "48.48"
$ws.Range("E9").Value = "  +5.33%  "
$ws.Range("D10").Value = "0.08232"
$ws.Range("E10").Value = "  +5.22%  "
$ws.Range("D11").Value = "1.039"
$ws.Range("E11").Value = "  +8.36%  "
$ws.Range("D12").Value = "22.60"
$ws.Range("E12").Value = "  +8.12%  "
$ws.Range("D13").Value = "1.927.61"
$ws.Range("E13").Value = "  +6.65%  "
$ws.Range("D14").Value = "6.175"
$ws.Range("E14").Value = "  +5.87%  "
$ws.Range("D15").Value = "7.403"
$ws.Range("D16").Value = "92.32"
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("D18").Value = "0.00001062"
$ws.Range("E18").Value = "  +4.33%  "
$ws.Range("D19").Value = "0.06665"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "18.02"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "29.615.86"
$ws.Range("E22").Value = "  +8.33%  "
$ws.Range("D23").Value = "5.608"
$ws.Range("E23").Value = "  +6.35%  "
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  +4.22%  "
$ws.Range("D25").Value = "2.285"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "2.180.88"
$ws.Range("E26").Value = "  +7.63%  "
$ws.Range("D27").Value = "160.59"
$ws.Range("E27").Value = "  +3.18%  "
$ws.Range("D28").Value = "20.15"
$ws.Range("E28").Value = "  +4.66%  "
$ws.Range("D29").Value = "2.189"
$ws.Range("E29").Value = "  +7.59%  "
$ws.Range("D30").Value = "5.619"
$ws.Range("E30").Value = "  +7.55%  "
$ws.Range("D31").Value = "122.12"
$ws.Range("E31").Value = "  +3.98%  "
$ws.Range("D32").Value = "1.023"
$ws.Range("E32").Value = "  +10.26%  "
$ws.Range("D33").Value = "0.09647"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("D34").Value = "1.463"
$ws.Range("E34").Value = "  +11.96%  "
$ws.Range("D35").Value = "3.686"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("D36").Value = "5.485"
$ws.Range("E36").Value = "  +5.41%  "
$ws.Range("D37").Value = "0.06282"
$ws.Range("E37").Value = "  +6.59%  "
$ws.Range("D38").Value = "0.02316"
$ws.Range("E38").Value = "  +6.95%  "
$ws.Range("D39").Value = "8.584"
$ws.Range("E39").Value = "  +6.90%  "
$ws.Range("E40").Value = "  +5.01%  "
$ws.Range("D41").Value = "0.6097"
$ws.Range("E41").Value = "  +6.65%  "
$ws.Range("D42").Value = "10.68"
$ws.Range("E42").Value = "  +8.14%  "
$ws.Range("D43").Value = "0.1904"
$ws.Range("E43").Value = "  +5.11%  "
$ws.Range("D45").Value = "1.274"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.62"
$ws.Range("E46").Value = "  +7.04%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5708"
$ws.Range("E47").Value = "  +6.37%  "
$ws.Range("D48").Value = "2.329"
$ws.Range("E48").Value = "  +31.16%  "
$ws.Range("D49").Value = "1.996"
$ws.Range("E49").Value = "  +7.36%  "
$ws.Range("D50").Value = "0.07405"
$ws.Range("E50").Value = "  +12.75%  "
$ws.Range("D51").Value = "114.24"
$ws.Range("E51").Value = "  +4.46%  "
